# Add a new "2020" column (N) to the SDG 6.3.1 indicator table, mirroring
# the formatting already used in the adjacent "2019" column (M) and filling
# in the new year's values for the Kyrgyz Republic and each oblast / city.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (thin bottom-border spacer row above the header) ---------------
$ws.Range("M3").Copy($ws.Range("N3"))

# --- Row 4 (year header row) ----------------------------------------------
$ws.Range("M4").Copy($ws.Range("N4"))
$ws.Range("N4").Value = 2020

# --- Row 5 (Kyrgyz Republic) -----------------------------------------------
$ws.Range("M5").Copy($ws.Range("N5"))
$ws.Range("N5").Value = 95.9

# --- Row 6 (Batken oblast) --------------------------------------------------
$ws.Range("M6").Copy($ws.Range("N6"))
$ws.Range("N6").Value = 96.2
$ws.Range("N6").NumberFormat = "0.0"

# --- Row 7 (Djalal-Abad oblast) ---------------------------------------------
$ws.Range("M7").Copy($ws.Range("N7"))
$ws.Range("N7").Value = 62.3
$ws.Range("N7").NumberFormat = "0.0"

# --- Row 8 (Ysyk-Kul oblast) -------------------------------------------------
$ws.Range("M8").Copy($ws.Range("N8"))
$ws.Range("N8").Value = 100
$ws.Range("N8").NumberFormat = "0.0"

# --- Row 9 (Naryn oblast) ----------------------------------------------------
$ws.Range("M9").Copy($ws.Range("N9"))
$ws.Range("N9").Value = 100
$ws.Range("N9").NumberFormat = "0.0"

# --- Row 10 (Osh oblast) -- data not available, keep the "-" placeholder ----
$ws.Range("M10").Copy($ws.Range("N10"))
$ws.Range("N10").NumberFormat = "0.0"

# --- Row 11 (Talas oblast) ---------------------------------------------------
$ws.Range("M11").Copy($ws.Range("N11"))
$ws.Range("N11").Value = 100
$ws.Range("N11").NumberFormat = "0.0"

# --- Row 12 (Chui oblast) -----------------------------------------------------
$ws.Range("M12").Copy($ws.Range("N12"))
$ws.Range("N12").Value = 62.7
$ws.Range("N12").NumberFormat = "0.0"

# --- Row 13 (Bishkek city) -----------------------------------------------------
$ws.Range("M13").Copy($ws.Range("N13"))
$ws.Range("N13").Value = 100
$ws.Range("N13").NumberFormat = "0.0"

# --- Row 14 (Osh city) -- data not available, keep the "-" placeholder ------
$ws.Range("M14").Copy($ws.Range("N14"))
$ws.Range("N14").NumberFormat = "0.0"

# Leave the selection on the newly-entered cell, like the original author did.
$ws.Range("N3").Select()
